# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Price ("D") and 1h-volume ("E") columns are plain text in this sheet, so
# any value that looks like a bare number (e.g. "224.13") is written with a
# leading apostrophe to force Excel to keep it as text instead of coercing
# it to a numeric cell (standard Excel "quote prefix" behaviour).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.026.31"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "1.782.10"
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'224.13"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'32.46"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.283"
$ws.Range("E9").Value = "  -3.84%  "
$ws.Range("D10").Value = "'0.0703"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "2.037.45"
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("D13").Value = "1.777.51"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").Value = "'10.79"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.621"
$ws.Range("E15").Value = "  -4.77%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "33.999.66"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "'4.14"
$ws.Range("E17").Value = "  -5.12%  "
$ws.Range("D18").Value = "'67.53"
$ws.Range("E18").Value = "  -3.45%  "
$ws.Range("D19").Value = "'242.52"
$ws.Range("E19").Value = "  -4.22%  "
$ws.Range("E20").Value = "  -3.44%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'10.66"
$ws.Range("E22").Value = "  -5.69%  "
$ws.Range("D23").Value = "'4.08"
$ws.Range("E23").Value = "  -5.80%  "
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "'159.34"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").Value = "'16.26"
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("D27").Value = "'7.02"
$ws.Range("E27").Value = "  -3.47%  "
$ws.Range("E28").Value = "  -2.84%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "'0.0512"
$ws.Range("E30").Value = "  -5.01%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("E33").Value = "  -4.49%  "
$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  -7.50%  "
$ws.Range("D35").Value = "1.390.46"
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("D38").Value = "'0.0185"
$ws.Range("E38").Value = "  -5.04%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "'2.35"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.20"
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.910"
$ws.Range("E42").Value = "  -6.92%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'78.24"
$ws.Range("E43").Value = "  -5.85%  "
$ws.Range("D44").Value = "0.0₆0147"
$ws.Range("E44").Value = "  +14.75%  "
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").Value = "'0.0497"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "'106.81"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "'5.85"
$ws.Range("E48").Value = "  -4.67%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'12.24"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "1.938.29"
$ws.Range("E50").Value = "  -3.54%  "
$ws.Range("E51").Value = "  -0.35%  "
